$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = '51.353.13'
$ws.Cells.Item(2, 5).Value = '  -1.12%  '

$ws.Cells.Item(3, 4).Value = '2.773.83'
$ws.Cells.Item(3, 5).Value = '  -0.64%  '

$ws.Cells.Item(4, 5).Value = '  -0.04%  '

$ws.Cells.Item(5, 4).Value = '''352.37'
$ws.Cells.Item(5, 5).Value = '  -2.57%  '

$ws.Cells.Item(6, 4).Value = '''108.64'
$ws.Cells.Item(6, 5).Value = '  -0.77%  '

$ws.Cells.Item(7, 5).Value = '  -2.05%  '

$ws.Cells.Item(8, 4).Value = '''0.999'
$ws.Cells.Item(8, 5).Value = '  -0.04%  '

$ws.Cells.Item(9, 4).Value = '''0.618'
$ws.Cells.Item(9, 5).Value = '  +4.62%  '

$ws.Cells.Item(10, 4).Value = '''39.17'
$ws.Cells.Item(10, 5).Value = '  -2.35%  '

$ws.Cells.Item(11, 5).Value = '  +1.51%  '

$ws.Cells.Item(12, 4).Value = '''0.0832'
$ws.Cells.Item(12, 5).Value = '  -1.92%  '

$ws.Cells.Item(13, 4).Value = '''19.85'
$ws.Cells.Item(13, 5).Value = '  +1.89%  '

$ws.Cells.Item(14, 4).Value = '''7.75'
$ws.Cells.Item(14, 5).Value = '  +2.49%  '

$ws.Cells.Item(15, 4).Value = '3.205.43'
$ws.Cells.Item(15, 5).Value = '  -0.73%  '

$ws.Cells.Item(16, 4).Value = '2.765.61'
$ws.Cells.Item(16, 5).Value = '  -1.12%  '

$ws.Cells.Item(17, 4).Value = '''0.922'
$ws.Cells.Item(17, 5).Value = '  -1.81%  '

$ws.Cells.Item(18, 4).Value = '51.263.61'
$ws.Cells.Item(18, 5).Value = '  -1.26%  '

$ws.Cells.Item(19, 4).Value = '''7.71'
$ws.Cells.Item(19, 5).Value = '  +3.13%  '

$ws.Cells.Item(20, 5).Value = '  +0.41%  '

$ws.Cells.Item(21, 4).Value = '''13.36'
$ws.Cells.Item(21, 5).Value = '  +1.86%  '

$ws.Cells.Item(22, 4).Value = '0.0₃0963'
$ws.Cells.Item(22, 5).Value = '  -1.33%  '

$ws.Cells.Item(23, 4).Value = '''70.24'
$ws.Cells.Item(23, 5).Value = '  -0.18%  '

$ws.Cells.Item(24, 4).Value = '''265.86'
$ws.Cells.Item(24, 5).Value = '  -1.53%  '

$ws.Cells.Item(25, 5).Value = '  -0.13%  '

$ws.Cells.Item(26, 5).Value = '  -0.13%  '

$ws.Cells.Item(27, 4).Value = '''25.78'
$ws.Cells.Item(27, 5).Value = '  -2.87%  '

$ws.Cells.Item(28, 4).Value = '''0.163'
$ws.Cells.Item(28, 5).Value = '  +0.48%  '

$ws.Cells.Item(29, 4).Value = '''10.24'
$ws.Cells.Item(29, 5).Value = '  -0.35%  '

$ws.Cells.Item(30, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(30, 4).Value = '''36.93'
$ws.Cells.Item(30, 5).Value = '  +7.91%  '

$ws.Cells.Item(31, 2).Value = 'Toncoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(31, 4).Value = '''2.23'
$ws.Cells.Item(31, 5).Value = '  -1.80%  '

$ws.Cells.Item(32, 5).Value = '  +8.96%  '

$ws.Cells.Item(33, 4).Value = '''51.83'
$ws.Cells.Item(33, 5).Value = '  -0.35%  '

$ws.Cells.Item(34, 5).Value = '  -5.97%  '

$ws.Cells.Item(35, 4).Value = '''5.54'
$ws.Cells.Item(35, 5).Value = '  +5.60%  '

$ws.Cells.Item(36, 5).Value = '  +0.00%  '

$ws.Cells.Item(37, 4).Value = '''0.0836'

$ws.Cells.Item(38, 4).Value = '''18.34'
$ws.Cells.Item(38, 5).Value = '  -3.55%  '

$ws.Cells.Item(39, 5).Value = '  -3.75%  '

$ws.Cells.Item(40, 5).Value = '  -2.10%  '

$ws.Cells.Item(41, 5).Value = '  -1.37%  '

$ws.Cells.Item(42, 5).Value = '  -4.44%  '

$ws.Cells.Item(43, 4).Value = '''119.95'
$ws.Cells.Item(43, 5).Value = '  +0.17%  '

$ws.Cells.Item(44, 5).Value = '  -2.69%  '

$ws.Cells.Item(45, 4).Value = '''21.74'
$ws.Cells.Item(45, 5).Value = '  -1.18%  '

$ws.Cells.Item(46, 4).Value = '2.128.01'
$ws.Cells.Item(46, 5).Value = '  +2.12%  '

$ws.Cells.Item(47, 2).Value = 'NEARProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(47, 4).Value = '''3.35'
$ws.Cells.Item(47, 5).Value = '  +2.98%  '

$ws.Cells.Item(48, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(48, 4).Value = '''2.34'
$ws.Cells.Item(48, 5).Value = '  +5.83%  '

$ws.Cells.Item(49, 4).Value = '''0.225'
$ws.Cells.Item(49, 5).Value = '  +18.42%  '

$ws.Cells.Item(50, 5).Value = '  -6.06%  '

$ws.Cells.Item(51, 5).Value = '  +8.82%  '
